# Applies the "Update finish generate form A dan form B" edit:
#  - Kepala [onshow. nama_satker_dua;noerr] ... Kepala Bidang Penyelenggaraan
#        -> [onshow. jabatan_kepala_satker;noerr] ... [onshow. jabatan_kepala_bidang;noerr]
#  - Drop one separating tab between nama_kepala_satker and nama_kepala_bidang fields
#  - Prefix both nip_kepala_satker / nip_kepala_bidang fields with "NIP " (consuming one tab)
#  - Refresh the "Lembar Ke-2 untuk Pusdiklat" paragraph so the stale lastRenderedPageBreak
#    cache marker is dropped

$d = $word.ActiveDocument

# 1) "Kepala [onshow. nama_satker_dua" -> "[onshow. jabatan_kepala_satker"
$r1 = $d.Content
$found1 = $r1.Find.Execute("Kepala [onshow. nama_satker_dua", $false, $false, $false, $false, $false, `
                            $true, 1, $false, "[onshow. jabatan_kepala_satker", 2)
Write-Host "jabatan_kepala_satker replace: $found1"

# 2) "Kepala Bidang Penyelenggaraan" -> "[onshow. jabatan_kepala_bidang;noerr]"
$r2 = $d.Content
$found2 = $r2.Find.Execute("Kepala Bidang Penyelenggaraan", $false, $false, $false, $false, $false, `
                            $true, 1, $false, "[onshow. jabatan_kepala_bidang;noerr]", 2)
Write-Host "jabatan_kepala_bidang replace: $found2"

# 3) Remove one of the five tabs between nama_kepala_satker and nama_kepala_bidang fields
$r3 = $d.Content
$found3 = $r3.Find.Execute("nama_kepala_satker;noerr]", $false, $false, $false, $false, $false, $false)
Write-Host "nama_kepala_satker find: $found3"
$tab1 = $d.Range($r3.End, $r3.End + 1)
Write-Host "tab1 text is tab: $($tab1.Text -eq [char]9)"
$tab1.Delete()

# 4) Prefix nip_kepala_satker field with "NIP "
$r4 = $d.Content
$found4 = $r4.Find.Execute("[onshow. nip_kepala_satker", $false, $false, $false, $false, $false, $false)
Write-Host "nip_kepala_satker find: $found4"
$ins4 = $d.Range($r4.Start, $r4.Start)
$ins4.InsertBefore("NIP ")

# 5) Replace the tab immediately before nip_kepala_bidang field with "NIP "
$r5 = $d.Content
$found5 = $r5.Find.Execute("[onshow. nip_kepala_bidang", $false, $false, $false, $false, $false, $false)
Write-Host "nip_kepala_bidang find: $found5"
$tab2 = $d.Range($r5.Start - 1, $r5.Start)
Write-Host "tab2 text is tab: $($tab2.Text -eq [char]9)"
$tab2.Text = "NIP "

# 6) Touch the "Lembar Ke-2 untuk Pusdiklat" paragraph so the stale lastRenderedPageBreak
#    render-cache marker gets dropped when the document is re-saved.
$r6 = $d.Content
$found6 = $r6.Find.Execute("Lembar Ke-2 untuk Pusdiklat", $false, $false, $false, $false, $false, $false)
Write-Host "Lembar Ke-2 find: $found6"
$r6.InsertAfter("")

Write-Host "Done"
